$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 702 (shifts existing rows
# 702..732 down to 703..733).
$ws.Rows.Item(702).Insert()

# Populate the newly inserted row with the new "Asterix" / "1a (guarda)"
# price observation (dated 45147) for Papa @ Feria Lagunitas de Puerto Montt.
$ws.Cells.Item(702, 1).Value = 4
$ws.Cells.Item(702, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(702, 3).Value = "Los Lagos"
$ws.Cells.Item(702, 4).Value = 45147
$ws.Cells.Item(702, 5).Value = 10
$ws.Cells.Item(702, 6).Value = 100114001
$ws.Cells.Item(702, 7).Value = "Papa"
$ws.Cells.Item(702, 8).Value = "Asterix"
$ws.Cells.Item(702, 9).Value = "1a (guarda)"
$ws.Cells.Item(702, 10).Value = 150
$ws.Cells.Item(702, 11).Value = 19000
$ws.Cells.Item(702, 12).Value = 20000
$ws.Cells.Item(702, 13).Value = 19467
$ws.Cells.Item(702, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(702, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(702, 16).Value = 779
$ws.Cells.Item(702, 17).Value = 25
$ws.Cells.Item(702, 18).Value = "Hortaliza"
